$wb = $excel.ActiveWorkbook
$ws5 = $wb.Worksheets.Item(5)
$ws3 = $wb.Worksheets.Item(3)

# --- Row 1 (header): rename existing columns, add new trailing columns ---
$ws5.Range("B1").Value = "species"
$ws5.Range("C1").Value = "debtor"
$ws5.Range("D1").Value = "owner"
$ws5.Range("E1").Value = "total"
$ws5.Range("F1").Value = "register_date"
$ws5.Range("G1").Value = "register_reason"
$ws3.Range("H1").Copy($ws5.Range("H1"))
$ws3.Range("I1").Copy($ws5.Range("I1"))
$ws3.Range("J1").Copy($ws5.Range("J1"))
$ws3.Range("K1").Copy($ws5.Range("K1"))
$ws3.Range("L1").Copy($ws5.Range("L1"))
$ws3.Range("M1").Copy($ws5.Range("M1"))
$ws3.Range("N1").Copy($ws5.Range("N1"))

# --- Row 2 (existing loan row) ---
$ws5.Range("B2").Value = "房屋貸款"
$ws5.Range("D2").Value = "台新銀行南京東路分行"
$ws5.Range("F2").Value = "100年02月17日"
$ws5.Range("H2").Value = "debt"
$ws5.Range("I2").Value = "normal"
$ws3.Range("J2").Copy($ws5.Range("J2"))
$ws5.Range("K2").Value = "李慶華"
$ws5.Range("L2").Value = 607
$ws5.Range("M2").Value = "tmp68f81"
$ws5.Range("N2").Value = 99

# --- Row 3 (existing parking-space loan row) ---
$ws5.Range("D3").Value = "台新銀行南京東路分行"
$ws5.Range("F3").Value = "100年04月14曰"
$ws5.Range("H3").Value = "debt"
$ws5.Range("I3").Value = "normal"
$ws3.Range("J3").Copy($ws5.Range("J3"))
$ws5.Range("K3").Value = "李慶華"
$ws5.Range("L3").Value = 607
$ws5.Range("M3").Value = "tmp68f81"
$ws5.Range("N3").Value = 100
